# Sprint Review Protocol #2 - fill in sprint metadata + mark reviewed items
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint metadata (column C, rows 4-12) ---
$ws.Range("C4").Value = "#2"
$ws.Range("C5").Value = "25.11.2021, 16:53:00"
$ws.Range("C6").Value = 29
$ws.Range("C7").Value = "Stefan Düx"
$ws.Range("C8").Value = "Jessica Isabella Görög"
$ws.Range("C9").Value = "Dominic Grabner"
$ws.Range("C10").Value = "Rebekka Tscheppen"
$ws.Range("C11").Value = "Lukas Varga"
$ws.Range("C12").Value = "Lukas Rohatsch"

# --- Mark requirement rows as reviewed (column F), except row 20 ---
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("F27").Value = 1

# --- Update the view: scroll near top and select I12 ---
$ws.Activate()
$ws.Range("I12").Select()
